# "checkout stripe finished, admin dashbord"
# Adds the new "Raw lamp" product row (row 4) to the products sheet,
# matches column formatting of the existing "Raw table" row (row 3),
# and leaves the sheet scrolled/selected near the new row as the author did.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new product row -------------------------------------------------
# Name / Desc / Price first, then Sku before Img_path -- this mirrors the
# order the strings were actually typed in (shared-string insertion order).
$ws.Range("A4").Value = "Raw lamp"
$ws.Range("B4").Value = "Brutalist pure lamp."
$ws.Range("C4").Value = 340
$ws.Range("E4").Value = "CC-FA-260212-01"
$ws.Range("D4").Value = "assets/products/prod_3_1.png;assets/products/prod_3_2.png;assets/products/prod_3_3.png;assets/products/prod_3_4.png"
$ws.Range("F4").Value = "Two-piece reversible layout;Monolithic concrete volume;Chamfered facets for sharp light;Micro-texture sealed top;Anti-slip pads + levelers"

# Carry over the same cell formatting used on row 3 for the img_path /
# features columns (bold sku font + wrap text respectively).
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D4").PasteSpecial(-4122) | Out-Null
$ws.Range("F3").Copy() | Out-Null
$ws.Range("F4").PasteSpecial(-4122) | Out-Null

# Row 3 is a wrapped 2-line row (height 29); match it on the new row too.
$ws.Rows.Item(4).RowHeight = $ws.Rows.Item(3).RowHeight

$excel.CutCopyMode = $false

# --- view state: scroll right a bit and land the selection on F4 -----
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
$ws.Range("F4").Select() | Out-Null

# --- page setup (print dialog touched paper size / orientation) ------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
